$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.236.80"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "2.064.58"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.22"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.59"
$ws.Range("E8").Value = "  -2.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "2.368.79"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.60"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.70"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.776"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.13"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "2.067.30"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "37.207.91"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("E19").Value = "  +2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.44"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.48"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.40"
$ws.Range("E26").Value = "  +2.18%  "

$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.42"
$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("E29").Value = "  -3.13%  "

$ws.Range("E30").Value = "  -2.48%  "

$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.48"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  +3.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0614"
$ws.Range("E34").Value = "  -2.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  -4.65%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.64"
$ws.Range("E39").Value = "  -4.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "1.477.83"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.92"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.32"
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("E44").Value = "  +3.06%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("E45").Value = "  -3.42%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0211"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.16"
$ws.Range("E48").Value = "  -6.15%  "

$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.68"
$ws.Range("E51").Value = "  -0.66%  "
